# ---------------------------------------------------------------------------
# Applies the "ActiviteEnseignement" StructureDefinition refresh:
#   1. Metadata sheet: bump the "Date" value to the new generation timestamp.
#   2. Elements sheet: append a new row describing the
#      "ActiviteEnseignement.EntiteGeographique" element (a reference to the
#      EntiteGeographique class), and widen the Type(s) column so the new
#      (long) URL fits.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Metadata!B8 ("Date" row) -------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-10-30T16:36:55+00:00"

# --- 2. Elements: append row 10 ---------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

$lastRow = 9
$newRow = $lastRow + 1

# Duplicate row 9 (values first, then formats) so the new row inherits the
# exact same "shape" (which columns carry text vs. stay blank) and the s="2"
# body-row style used throughout the table.
$srcRow = "A" + $lastRow + ":AJ" + $lastRow
$dstRow = "A" + $newRow + ":AJ" + $newRow

$ws.Range($srcRow).Copy()
$ws.Range($dstRow).PasteSpecial(-4104)   # xlPasteValues
$ws.Range($srcRow).Copy()
$ws.Range($dstRow).PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Overwrite with the new element's data.
$ws.Range("A" + $newRow).Value = "ActiviteEnseignement.EntiteGeographique"
$ws.Range("B" + $newRow).Value = "ActiviteEnseignement.EntiteGeographique"

# Min / Max / Base Min / Base Max are "1".."1" — force text storage (these
# columns are shared-string text everywhere else in the sheet) by briefly
# marking the cells as Text before writing the numeric-looking string.
$ws.Range("F" + $newRow + ":G" + $newRow).NumberFormat = "@"
$ws.Range("F" + $newRow).Value = "1"
$ws.Range("G" + $newRow).Value = "1"

$ws.Range("K" + $newRow).Value = "https://interop.esante.gouv.fr/ig/mos/StructureDefinition/EntiteGeographique`n"
$ws.Range("L" + $newRow).Value = "Lien vers la classe EntiteGeographique"
$ws.Range("M" + $newRow).Value = "Lien vers la classe EntiteGeographique"

$ws.Range("AF" + $newRow).Value = "ActiviteEnseignement.EntiteGeographique"
$ws.Range("AG" + $newRow + ":AH" + $newRow).NumberFormat = "@"
$ws.Range("AG" + $newRow).Value = "1"
$ws.Range("AH" + $newRow).Value = "1"

# Re-apply the body-row format on top so every cell (including the ones we
# just touched) ends up back on the shared s="2" style instead of the
# transient Text-number-format variant.
$ws.Range($srcRow).Copy()
$ws.Range($dstRow).PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Let the row height go back to "auto" rather than sticking at whatever the
# paste operations left behind.
$ws.Rows.Item($newRow).AutoFit()

# Column K ("Type(s)") now holds a long URL — widen it to fit, matching the
# width the real workbook ends up with once the column is re-fitted.
# (61.65625 is the target stored width; this engine re-quantizes whatever
# we request at save time, so 60.75 is the input that lands closest to it.)
$ws.Columns.Item(11).ColumnWidth = 60.75
